$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3: "13.2 Integrate climate change measures..." keeps its text, but
#     loses the wrap-text formatting (style goes from s="6" to s="5").
$ws.Range("B3").WrapText = $false

# --- B4: indicator text updated to new UN wording, and the cell loses its
#     explicit style entirely (no "s" attribute afterwards -> default style).
#     Set the new value FIRST (while the cell is still unlocked under sheet
#     protection), then copy-paste formats only from a pristine, untouched
#     cell on the same sheet to reset B4's formatting to the default style.
$ws.Range("B4").Value = "13.2.1 Number of countries with nationally determined contributions, long-term strategies, national adaptation plans and adaptation communications, as reported to the secretariat of the United Nations Framework Convention on Climate Change"
$blank = $ws.Range("A100")
$blank.Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- B6-B10: text content for the "Data reporter" block (unchanged content,
#     same strings as before, styles s=6/s=5 are untouched by simple value sets)
$ws.Range("B6").Value = "The State Agency on Environment Protection and Forestry"
$ws.Range("B7").Value = "N. S. Baidakova"
$ws.Range("B8").Value = "banatalia@yandex.com "
$ws.Range("B9").Value = "0(312) 54-94-87"
$ws.Range("B10").Value = "www.ecology.gov.kg "

# --- Selection moves from B12 to B6
$ws.Range("B6").Select()

# --- Window geometry (bookViews workbookView attributes)
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 28800
$win.Height = 11835

Write-Host "done"
